# Report details.xlsx update: add 4 new test-run rows (13-16) for
# 2022-01-14 and 2022-01-18, and move the view/selection down to them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: Development run on 2022-01-14 (serial 44575), 2 fail cases ---
$ws.Range("A13").Value = "1/14/2022"
$ws.Range("B13").Value = "Development"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "After execution all test cases pass"
$ws.Range("E13").Value = "Test cases initially fail because of page load affected by network"

# --- Row 14: Production run on 2022-01-14 (serial 44575), 5 fail cases ---
$ws.Range("A14").Value = "1/14/2022"
$ws.Range("B14").Value = "Production"
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = "After execution all test cases pass"
$ws.Range("E14").Value = "Test cases initially fail because of page load affected by network"

# --- Row 15: Development run on 2022-01-18 (serial 44579), 2 fail cases ---
$ws.Range("A15").Value = "1/18/2022"
$ws.Range("B15").Value = "Development"
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "After execution all test cases pass"
$ws.Range("E15").Value = "Test cases initially fail because of page load affected by network"

# --- Row 16: Production run on 2022-01-18 (serial 44579), 2 fail cases ---
$ws.Range("A16").Value = "1/18/2022"
$ws.Range("B16").Value = "Production"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = "After execution all test cases pass"
$ws.Range("E16").Value = "Test cases initially fail because of page load affected by network"

# Copy the formatting (date/number formats, wrap-text, borders/fill) from the
# previous data row (12), which already has the same A:E layout, onto the
# four new rows.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The new rows hold wrapped comment text, same as the other data rows, so
# they use the taller 75pt row height.
$ws.Rows.Item(13).RowHeight = 75
$ws.Rows.Item(14).RowHeight = 75
$ws.Rows.Item(15).RowHeight = 75
$ws.Rows.Item(16).RowHeight = 75

# Scroll the view down to the newly added rows and select D15:E16.
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("D15:E16").Select()

# Widen the saved window.
$excel.ActiveWindow.Width = 13500
